# Apply the changes described in the diff:
# 1. Set cell C530 to "Age" (new value added in row 530).
# 2. Delete the blank placeholder row 642 entirely, which shifts all the
#    following rows (previously 643-665) up by one (becoming 642-664),
#    and reduces the used range from A1:O665 to A1:O664.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Add the "Age" label in column C of row 530.
$ws.Range("C530").Value = "Age"

# 2. Remove the empty row 642 and shift everything below it up.
$ws.Rows("642").Delete()
